# "Literature Summary Table" restructure:
#  - drop the "Geography" column (old C)
#  - drop the "Result" column (old E), replaced with a new "Methodology" column
#  - the old "Signal" column (old D) slides left into C
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Geography" (column C) - Signal + Result shift left automatically
$ws.Range("C1").EntireColumn.Delete()

# Remove "Result" (now column D after the previous delete)
$ws.Range("D1").EntireColumn.Delete()

# New "Methodology" column header + the single populated note (Lemmon, 2015 row)
$ws.Range("D1").Value = "Methodology"
$ws.Range("D9").Value = "Grouped high yield stock and tests for yield effect at portfolio level.  Fama-MacBeth methodology and tests for yield effect after controlling for known factors. "

# Column widths: Signal keeps its old width, Methodology gets a wider column
$ws.Columns.Item(3).ColumnWidth = 17.17
$ws.Columns.Item(4).ColumnWidth = 37.83

# Refresh the sort range/state now that the table is only 4 columns wide
$rng = $ws.Range("A2:D11")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A11"))
$ws.Sort.SetRange($rng)
$ws.Sort.Apply()

# Selection moved onto the new note cell
[void]$ws.Range("D7").Select()
